# Đổi tên các sheet cho dễ theo dõi
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# Rename sheets
$ws1.Name = "Version Control"
$ws2.Name = "Project Plan"
$ws3.Name = "Iteration 1"

# Update selection/view state on each sheet (visiting each sheet clears any
# stale topLeftCell scroll position from the previous edit session)
$ws1.Activate()
$ws1.Range("C12").Select()

$ws3.Activate()
$ws3.Range("C12").Select()

# Leave "Project Plan" as the active sheet/tab
$ws2.Activate()
$ws2.Range("C9").Select()
